$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.06694645484776629
$ws.Range("D2").Value = 0.3380945829319444
$ws.Range("E2").Value = 0.04560622667644765
$ws.Range("F2").Value = 4.950746586021864
$ws.Range("G2").Value = 0.002595379632432633
$ws.Range("I2").Value = 4.172259820914093
$ws.Range("L2").Value = 0.06423515957794024
$ws.Range("M2").Value = 3.217101330326045
$ws.Range("N2").Value = 1.413707751416609

$ws.Range("C3").Value = 0.0592750218616942
$ws.Range("D3").Value = 0.3323401450663681
$ws.Range("E3").Value = 0.0425945261274876
$ws.Range("F3").Value = 4.786912834133659
$ws.Range("G3").Value = 0.002608374373772103
$ws.Range("I3").Value = 4.058759093110382
$ws.Range("L3").Value = 0.06224955109791352
$ws.Range("M3").Value = 2.983587420919974
$ws.Range("N3").Value = 1.342186129967388

$ws.Range("C4").Value = 0.0545575295000873
$ws.Range("D4").Value = 0.3291448820554876
$ws.Range("E4").Value = 0.04073627938670654
$ws.Range("F4").Value = 4.691117266246266
$ws.Range("G4").Value = 0.002616737082757062
$ws.Range("I4").Value = 3.993081151098693
$ws.Range("L4").Value = 0.06103312933895211
$ws.Range("M4").Value = 2.841585365850875
$ws.Range("N4").Value = 1.298753576954596

$ws.Range("C5").Value = 0.05263274539684915
$ws.Range("D5").Value = 0.3279261239148497
$ws.Range("E5").Value = 0.03997668394091747
$ws.Range("F5").Value = 4.653259618105096
$ws.Range("G5").Value = 0.002620242048736728
$ws.Range("I5").Value = 3.967304793846353
$ws.Range("L5").Value = 0.06053815395539885
$ws.Range("M5").Value = 2.784054939707488
$ws.Range("N5").Value = 1.281177600569208

$ws.Range("C6").Value = 0.05231297746340147
$ws.Range("D6").Value = 0.3277287359880177
$ws.Range("E6").Value = 0.03985041011627999
$ws.Range("F6").Value = 4.647043869426
$ws.Range("G6").Value = 0.002620829926502153
$ws.Range("I6").Value = 3.963083717934353
$ws.Range("L6").Value = 0.06045600857929401
$ws.Range("M6").Value = 2.774522112480071
$ws.Range("N6").Value = 1.278266618913847

$ws.Range("C7").Value = 0.05453158150083937
$ws.Range("D7").Value = 0.3291281101753185
$ws.Range("E7").Value = 0.04072604480015229
$ws.Range("F7").Value = 4.690601962742818
$ws.Range("G7").Value = 0.00261678395817877
$ws.Range("I7").Value = 3.992729549205791
$ws.Range("L7").Value = 0.061026450940318
$ws.Range("M7").Value = 2.840808140000775
$ws.Range("N7").Value = 1.298516039969911

$ws.Range("C8").Value = 0.0643024979012381
$ws.Range("D8").Value = 0.3360392937708383
$ws.Range("E8").Value = 0.04456962748237459
$ws.Range("F8").Value = 4.893244184678423
$ws.Range("G8").Value = 0.002599780885489844
$ws.Range("I8").Value = 4.132279290639303
$ws.Range("L8").Value = 0.06354997112855187
$ws.Range("M8").Value = 3.136294114788797
$ws.Range("N8").Value = 1.388948167637892

$ws.Range("C9").Value = 0.08343264707106357
$ws.Range("D9").Value = 0.3523514522009066
$ws.Range("E9").Value = 0.05203866304854188
$ws.Range("F9").Value = 5.329974994012503
$ws.Range("G9").Value = 0.00256945789839395
$ws.Range("I9").Value = 4.438755481021985
$ws.Range("L9").Value = 0.06851915543617082
$ws.Range("M9").Value = 3.727142568257221
$ws.Range("N9").Value = 1.570038496100807

$ws.Range("C10").Value = 0.09750786653569321
$ws.Range("D10").Value = 0.3661303760228236
$ws.Range("E10").Value = 0.05749026452753725
$ws.Range("F10").Value = 5.676678987301273
$ws.Range("G10").Value = 0.002548983845367392
$ws.Range("I10").Value = 4.685363911683964
$ws.Range("L10").Value = 0.07218129541893603
$ws.Range("M10").Value = 4.168927770780869
$ws.Range("N10").Value = 1.705297953045658

$ws.Range("C11").Value = 0.1039243247435593
$ws.Range("D11").Value = 0.3728139259267209
$ws.Range("E11").Value = 0.05996387984879448
$ws.Range("F11").Value = 5.840434975065818
$ws.Range("G11").Value = 0.00254005335619057
$ws.Range("I11").Value = 4.802534172299119
$ws.Range("L11").Value = 0.07384952009905987
$ws.Range("M11").Value = 4.371749322499397
$ws.Range("N11").Value = 1.76729793222313

$ws.Range("C12").Value = 0.1063567547499815
$ws.Range("D12").Value = 0.3754067132685464
$ws.Range("E12").Value = 0.06089976547637477
$ws.Range("F12").Value = 5.903348795310421
$ws.Range("G12").Value = 0.002536726074582002
$ws.Range("I12").Value = 4.847648031713106
$ws.Range("L12").Value = 0.07448153788870115
$ws.Range("M12").Value = 4.448832753353656
$ws.Range("N12").Value = 1.790841933237459

$ws.Range("C13").Value = 0.1058327577585487
$ws.Range("D13").Value = 0.3748455235473784
$ws.Range("E13").Value = 0.06069824063116513
$ws.Range("F13").Value = 5.889758452467163
$ws.Range("G13").Value = 0.002537440250192292
$ws.Range("I13").Value = 4.837898424572359
$ws.Range("L13").Value = 0.07434540884479901
$ws.Range("M13").Value = 4.432218819036166
$ws.Range("N13").Value = 1.785768402554083

$ws.Range("C14").Value = 0.1041243848442832
$ws.Range("D14").Value = 0.3730259835948857
$ws.Range("E14").Value = 0.06004089157711334
$ws.Range("F14").Value = 5.845592632752755
$ws.Range("G14").Value = 0.002539778530083367
$ws.Range("I14").Value = 4.806230647469164
$ws.Range("L14").Value = 0.07390151081074237
$ws.Range("M14").Value = 4.378085347484131
$ws.Range("N14").Value = 1.769233598384517

$ws.Range("C15").Value = 0.1030783244578686
$ws.Range("D15").Value = 0.3719195850374604
$ws.Range("E15").Value = 0.05963814253684063
$ws.Range("F15").Value = 5.818658438196849
$ws.Range("G15").Value = 0.002541217873274693
$ws.Range("I15").Value = 4.786930929999329
$ws.Range("L15").Value = 0.07362964819754581
$ws.Range("M15").Value = 4.344963807993764
$ws.Range("N15").Value = 1.759114096641667

$ws.Range("C16").Value = 0.09708886863907651
$ws.Range("D16").Value = 0.3657021305420756
$ws.Range("E16").Value = 0.05732848815120306
$ws.Range("F16").Value = 5.666101591765937
$ws.Range("G16").Value = 0.002549575133623728
$ws.Range("I16").Value = 4.677809186016731
$ws.Range("L16").Value = 0.07207231640371958
$ws.Range("M16").Value = 4.15571126220334
$ws.Range("N16").Value = 1.701255423044643

$ws.Range("C17").Value = 0.09341850051168876
$ws.Range("D17").Value = 0.3619958209788479
$ws.Range("E17").Value = 0.05591002979421589
$ws.Range("F17").Value = 5.574083775118311
$ws.Range("G17").Value = 0.00255479977285545
$ws.Range("I17").Value = 4.612162619410213
$ws.Range("L17").Value = 0.07111750858075538
$ws.Range("M17").Value = 4.040094218667633
$ws.Range("N17").Value = 1.665880156608722

$ws.Range("C18").Value = 0.09130863537332345
$ws.Range("D18").Value = 0.359903021468341
$ws.Range("E18").Value = 0.05509356688321532
$ws.Range("F18").Value = 5.521723793211549
$ws.Range("G18").Value = 0.002557840958676909
$ws.Range("I18").Value = 4.574872165794204
$ws.Range("L18").Value = 0.07056854667615653
$ws.Range("M18").Value = 3.973767777646344
$ws.Range("N18").Value = 1.645577546462732

$ws.Range("C19").Value = 0.09059446209825239
$ws.Range("D19").Value = 0.3592010653127318
$ws.Range("E19").Value = 0.05481702020256662
$ws.Range("F19").Value = 5.504091819711391
$ws.Range("G19").Value = 0.002558876872336297
$ws.Range("I19").Value = 4.562325785612046
$ws.Range("L19").Value = 0.07038271613193814
$ws.Range("M19").Value = 3.951340181726351
$ws.Range("N19").Value = 1.638711091475699

$ws.Range("C20").Value = 0.09380908477606908
$ws.Range("D20").Value = 0.3623863153528646
$ws.Range("E20").Value = 0.05606108887236161
$ws.Range("F20").Value = 5.583820379127701
$ws.Range("G20").Value = 0.002554239868140837
$ws.Range("I20").Value = 4.619102209962705
$ws.Range("L20").Value = 0.07121912715427925
$ws.Range("M20").Value = 4.052383808749369
$ws.Range("N20").Value = 1.669641340308033

$ws.Range("C21").Value = 0.1046260975313089
$ws.Range("D21").Value = 0.3735587300115242
$ws.Range("E21").Value = 0.06023399237920302
$ws.Range("F21").Value = 5.858540415332641
$ws.Range("G21").Value = 0.00253909024630938
$ws.Range("I21").Value = 4.815511828975417
$ws.Range("L21").Value = 0.07403188661488258
$ws.Range("M21").Value = 4.39397797210367
$ws.Range("N21").Value = 1.774088493501722

$ws.Range("C22").Value = 0.111711471808178
$ws.Range("D22").Value = 0.3812222341598499
$ws.Range("E22").Value = 0.06295650061222702
$ws.Range("F22").Value = 6.043366135454391
$ws.Range("G22").Value = 0.002529506473611432
$ws.Range("I22").Value = 4.948226848549837
$ws.Range("L22").Value = 0.07587190579486958
$ws.Range("M22").Value = 4.618864110441024
$ws.Range("N22").Value = 1.842734498500022

$ws.Range("C23").Value = 0.1079281855825514
$ws.Range("D23").Value = 0.377098258107452
$ws.Range("E23").Value = 0.06150384680079313
$ws.Range("F23").Value = 5.944226232374149
$ws.Range("G23").Value = 0.002534592678006686
$ws.Range("I23").Value = 4.876987157458075
$ws.Range("L23").Value = 0.07488970641890802
$ws.Range("M23").Value = 4.49868407344178
$ws.Range("N23").Value = 1.806062226132383

$ws.Range("C24").Value = 0.09363250081618446
$ws.Range("D24").Value = 0.3622096546227453
$ws.Range("E24").Value = 0.0559927980987176
$ws.Range("F24").Value = 5.579416773958087
$ws.Range("G24").Value = 0.002554492884237615
$ws.Range("I24").Value = 4.615963419480522
$ws.Range("L24").Value = 0.07117318548849738
$ws.Range("M24").Value = 4.046827239591693
$ws.Range("N24").Value = 1.667940799700972

$ws.Range("C25").Value = 0.07825676447131968
$ws.Range("D25").Value = 0.3476318194422561
$ws.Range("E25").Value = 0.05002482166145228
$ws.Range("F25").Value = 5.207427888473717
$ws.Range("G25").Value = 0.002577341534515871
$ws.Range("I25").Value = 4.352186985479904
$ws.Range("L25").Value = 0.06717279076266891
$ws.Range("M25").Value = 3.566007546668061
$ws.Range("N25").Value = 1.520655096504186
